{"js": "// Replace the text of each arithmetic expression in the 20x5 table,\n// in document order (row-major), while preserving the existing run\n// and paragraph formatting (font, size, alignment) of each cell.\nconst newValues = [\n  [\"45+47=\", \"85-58=\", \"5+79=\", \"71-23=\", \"85+8=\"],\n  [\"66-7=\", \"55-17=\", \"24-15=\", \"59+22=\", \"87-68=\"],\n  [\"61-35=\", \"73-25=\", \"5+57=\", \"67-48=\", \"18+46=\"],\n  [\"71-6=\", \"90-13=\", \"13+29=\", \"18+76=\", \"23+38=\"],\n  [\"67+17=\", \"39+43=\", \"57+34=\", \"27+29=\", \"68+17=\"],\n  [\"36+25=\", \"9+29=\", \"9+69=\", \"82-6=\", \"41-17=\"],\n  [\"80-2=\", \"12+19=\", \"32-14=\", \"58+4=\", \"26+18=\"],\n  [\"7+58=\", \"40-7=\", \"54+29=\", \"86-48=\", \"23-9=\"],\n  [\"71-35=\", \"83-27=\", \"51-28=\", \"68+28=\", \"50-12=\"],\n  [\"29+24=\", \"48+24=\", \"9+73=\", \"5+69=\", \"25+19=\"],\n  [\"19+42=\", \"80-25=\", \"62-13=\", \"19+2=\", \"56+6=\"],\n  [\"52-16=\", \"45+29=\", \"7+5=\", \"54+17=\", \"47+47=\"],\n  [\"91-87=\", \"27+6=\", \"22+19=\", \"82-75=\", \"51-19=\"],\n  [\"71-38=\", \"65-57=\", \"94-8=\", \"95-19=\", \"58+5=\"],\n  [\"96-49=\", \"28+33=\", \"37+58=\", \"18+17=\", \"77+19=\"],\n  [\"24+48=\", \"70-2=\", \"18+33=\", \"91-9=\", \"72-14=\"],\n  [\"31-16=\", \"84-36=\", \"29+32=\", \"62-54=\", \"40-8=\"],\n  [\"26+65=\", \"83-35=\", \"38+55=\", \"42+29=\", \"2+79=\"],\n  [\"36-17=\", \"60-22=\", \"13+18=\", \"41-33=\", \"25+16=\"],\n  [\"33+38=\", \"70-15=\", \"72-14=\", \"17+6=\", \"62+29=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const paragraph = cell.body.paragraphs.getFirst();\n    const range = paragraph.getRange();\n    // Replacing the paragraph's range (rather than the cell body) keeps\n    // the existing run's rPr (font/size) and the paragraph's pPr intact.\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the text of each arithmetic expression in the 20x5 table,\n# in document order (row-major), while preserving the existing run\n# and paragraph formatting (font, size, alignment) of each cell.\n$newValues = @(\n    @(\"45+47=\", \"85-58=\", \"5+79=\", \"71-23=\", \"85+8=\"),\n    @(\"66-7=\", \"55-17=\", \"24-15=\", \"59+22=\", \"87-68=\"),\n    @(\"61-35=\", \"73-25=\", \"5+57=\", \"67-48=\", \"18+46=\"),\n    @(\"71-6=\", \"90-13=\", \"13+29=\", \"18+76=\", \"23+38=\"),\n    @(\"67+17=\", \"39+43=\", \"57+34=\", \"27+29=\", \"68+17=\"),\n    @(\"36+25=\", \"9+29=\", \"9+69=\", \"82-6=\", \"41-17=\"),\n    @(\"80-2=\", \"12+19=\", \"32-14=\", \"58+4=\", \"26+18=\"),\n    @(\"7+58=\", \"40-7=\", \"54+29=\", \"86-48=\", \"23-9=\"),\n    @(\"71-35=\", \"83-27=\", \"51-28=\", \"68+28=\", \"50-12=\"),\n    @(\"29+24=\", \"48+24=\", \"9+73=\", \"5+69=\", \"25+19=\"),\n    @(\"19+42=\", \"80-25=\", \"62-13=\", \"19+2=\", \"56+6=\"),\n    @(\"52-16=\", \"45+29=\", \"7+5=\", \"54+17=\", \"47+47=\"),\n    @(\"91-87=\", \"27+6=\", \"22+19=\", \"82-75=\", \"51-19=\"),\n    @(\"71-38=\", \"65-57=\", \"94-8=\", \"95-19=\", \"58+5=\"),\n    @(\"96-49=\", \"28+33=\", \"37+58=\", \"18+17=\", \"77+19=\"),\n    @(\"24+48=\", \"70-2=\", \"18+33=\", \"91-9=\", \"72-14=\"),\n    @(\"31-16=\", \"84-36=\", \"29+32=\", \"62-54=\", \"40-8=\"),\n    @(\"26+65=\", \"83-35=\", \"38+55=\", \"42+29=\", \"2+79=\"),\n    @(\"36-17=\", \"60-22=\", \"13+18=\", \"41-33=\", \"25+16=\"),\n    @(\"33+38=\", \"70-15=\", \"72-14=\", \"17+6=\", \"62+29=\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n  $row = $newValues[$r - 1]\n  for ($c = 1; $c -le $row.Count; $c++) {\n    $cell = $tbl.Cell($r, $c)\n    # Assigning to Range.Text replaces only the text run content and\n    # keeps the existing rPr/pPr formatting on the cell's paragraph.\n    $cell.Range.Text = $row[$c - 1]\n  }\n}\n"}
